$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 31; this shifts current rows 31-67 down to 32-68
$ws.Rows.Item(31).Insert()

# Populate the newly inserted row 31 with the new weekly record.
# Columns A,B,C,E,F,G,H,I,N,O,Q,R carry the same values as the record that
# used to occupy row 31 (now shifted to row 32); only the date / volume /
# price columns (D,J,K,L,M,P) differ for this new entry.
$ws.Cells.Item(31, 1).Value = 2
$ws.Cells.Item(31, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(31, 3).Value = "Coquimbo"
$ws.Cells.Item(31, 4).Value = 44790
$ws.Cells.Item(31, 5).Value = 4
$ws.Cells.Item(31, 6).Value = 100112022
$ws.Cells.Item(31, 7).Value = "Arveja Verde"
$ws.Cells.Item(31, 8).Value = "Perfection"
$ws.Cells.Item(31, 9).Value = "Primera"
$ws.Cells.Item(31, 10).Value = 560
$ws.Cells.Item(31, 11).Value = 27000
$ws.Cells.Item(31, 12).Value = 29000
$ws.Cells.Item(31, 13).Value = 28000
$ws.Cells.Item(31, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(31, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(31, 16).Value = 1120
$ws.Cells.Item(31, 17).Value = 25
$ws.Cells.Item(31, 18).Value = "Hortaliza"
